# Stimulus update: swap the "face" image category for the "book" category,
# and spell out the abbreviated response-key codes used in column L
# (correct_ans) so they read as full position names instead of single
# letters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole

# "face//face_NN.jpg" -> "book//book_NN.jpg" everywhere on the sheet
# (promptFile / correctFile / dist_0XFile columns).
$ws.UsedRange.Replace("face", "book")

# correct_ans column: b -> center, y -> left, r -> right.
# Use whole-cell matching so we only touch cells that are exactly "b"/"y"/"r"
# and don't clobber substrings inside other filenames (e.g. "bug", "berry").
$col = $ws.Range("L2:L361")
$col.Replace("b", "center", $xlWhole)
$col.Replace("y", "left", $xlWhole)
$col.Replace("r", "right", $xlWhole)
